$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.890.24'
$ws.Range('E2').Value = '  -1.02%  '
$ws.Range('D3').Value = '1.896.62'
$ws.Range('E3').Value = '  -0.66%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7549'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '240.15'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.68%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3041'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.84%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.38'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.90%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06834'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.83%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07970'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('D12').Value = '1.907.74'
$ws.Range('E12').Value = '  -0.45%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7461'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.200'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.58%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.12'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.58%  '
$ws.Range('D16').Value = '29.894.61'
$ws.Range('E16').Value = '  -1.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.90'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.67%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.947'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '242.78'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.65%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007721'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.54%  '
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.939'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.36%  '
$ws.Range('E24').Value = '  -2.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '165.41'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.74'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.28%  '
$ws.Range('E27').Value = '  +2.48%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.024'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.404'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.514'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.277'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.018'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.43%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05349'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.27%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.247'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.69%  '
$ws.Range('E35').Value = '  -2.76%  '
$ws.Range('E36').Value = '  -1.62%  '
$ws.Range('E37').Value = '  -1.42%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.785'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.31%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.172'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.4398'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.82%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '72.21'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.28%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.001'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.06%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.904'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.83%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8245'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.29%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.09'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.22%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.548'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.62%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.751'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.93%  '
$ws.Range('D48').Value = '2.056.64'
$ws.Range('E48').Value = '  -1.64%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '36.16'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05967'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.56%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.463'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.28%  '
